# Apply the row-2 deletion + interest-count ("想去人数", column F) refresh
# to the "展览" (Exhibitions) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# Deltas to add to column F (想去人数) after the row has shifted up by one,
# keyed by the NEW row number (post-deletion) on each affected sheet.
$deltasSheet1 = @{
    3  = 22
    4  = 40
    5  = 40
    7  = 5
    16 = 1
    17 = 4
    22 = 3
    23 = -1
    24 = 19
    26 = 2
    28 = 2
    30 = 5
    31 = 2
    33 = 2
    35 = 2
    36 = 13
}

$deltasSheet4 = @{
    3  = 22
    4  = 40
    5  = 40
    7  = 5
    16 = 1
    17 = 4
    23 = 3
    24 = -1
    25 = 19
    27 = 2
    29 = 2
    32 = 5
    33 = 2
    35 = 2
    37 = 2
    38 = 13
}

function Apply-RowDeleteAndDeltas($sheetName, $deltas) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Rows.Item(2).Delete()

    foreach ($rowNum in $deltas.Keys) {
        $delta = $deltas[$rowNum]
        $cell = $ws.Cells.Item($rowNum, 6)
        $cell.Value2 = $cell.Value2 + $delta
    }
}

Apply-RowDeleteAndDeltas "展览" $deltasSheet1
Apply-RowDeleteAndDeltas "全部类型" $deltasSheet4
